$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 538.3333
$ws.Range("I33").Value = 538.3333
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 538.3333
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -309.3333
$ws.Range("N33").Value = $null
$ws.Range("H51").Value = 2938.6667
$ws.Range("I51").Value = 3266.6667
$ws.Range("J51").Value = 2856.6667
$ws.Range("K51").Value = 3266.6667
$ws.Range("L51").Value = 2856.6667
$ws.Range("M51").Value = -2782.6667
$ws.Range("N51").Value = -3824.6667
$ws.Range("H103").Value = 1312.5
$ws.Range("J103").Value = 1460
$ws.Range("L103").Value = 4380
$ws.Range("N103").Value = -5552
$ws.Range("H112").Value = 2209.3333
$ws.Range("I112").Value = 640
$ws.Range("J112").Value = 2321.4285
$ws.Range("K112").Value = 1920
$ws.Range("L112").Value = 6964.2855
$ws.Range("M112").Value = -812
$ws.Range("N112").Value = -9180.2855
$ws.Range("H134").Value = 63020
$ws.Range("J134").Value = 63020
$ws.Range("L134").Value = 63020
$ws.Range("N134").Value = -73160
$ws.Range("H137").Value = 27779512
$ws.Range("I137").Value = 1197.1428
$ws.Range("J137").Value = 125003620
$ws.Range("K137").Value = 3591.4284
$ws.Range("L137").Value = 375010860
$ws.Range("M137").Value = -1041.4284
$ws.Range("N137").Value = -375015960
$ws.Range("H138").Value = 2706628
$ws.Range("I138").Value = 5886005
$ws.Range("K138").Value = 17658015
$ws.Range("M138").Value = -17652875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4346.4224
$ws.Range("I32").Value = 4469.3896
$ws.Range("K32").Value = 4469.3896
$ws.Range("M32").Value = -4182.3896
$ws.Range("H102").Value = 2174
$ws.Range("I102").Value = 2174
$ws.Range("K102").Value = 2174
$ws.Range("M102").Value = -552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 9542
$ws.Range("J28").Value = 9542
$ws.Range("L28").Value = 9542
$ws.Range("N28").Value = -10130
$ws.Range("H98").Value = 59441.5
$ws.Range("J98").Value = 59441.5
$ws.Range("L98").Value = 59441.5
$ws.Range("N98").Value = -65431.5
$ws.Range("H107").Value = 1588.3636
$ws.Range("I107").Value = 1719.1111
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1719.1111
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 200.8888999999999
$ws.Range("N107").Value = -4840
$ws.Range("H134").Value = 19795.51
$ws.Range("I134").Value = 20118.5
$ws.Range("K134").Value = 60355.5
$ws.Range("M134").Value = -57820.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1786.3871
$ws.Range("I31").Value = 1493.6316
$ws.Range("J31").Value = 2249.9167
$ws.Range("K31").Value = 1493.6316
$ws.Range("L31").Value = 2249.9167
$ws.Range("M31").Value = -1198.6316
$ws.Range("N31").Value = -2839.9167
$ws.Range("H34").Value = 1786.3871
$ws.Range("I34").Value = 1493.6316
$ws.Range("J34").Value = 2249.9167
$ws.Range("K34").Value = 1493.6316
$ws.Range("L34").Value = 2249.9167
$ws.Range("M34").Value = -1291.6316
$ws.Range("N34").Value = -2653.9167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8517837
$ws.Range("I4").Value = 11227649
$ws.Range("J4").Value = 1285.7142
$ws.Range("K4").Value = 33682947
$ws.Range("L4").Value = 3857.1426
$ws.Range("M4").Value = -33682835
$ws.Range("N4").Value = -4081.1426
$ws.Range("H5").Value = 932.7241
$ws.Range("I5").Value = 783.2222
$ws.Range("K5").Value = 2349.6666
$ws.Range("M5").Value = -2237.6666
$ws.Range("H34").Value = 71429400
$ws.Range("I34").Value = 378.33334
$ws.Range("J34").Value = 90910040
$ws.Range("K34").Value = 1135.00002
$ws.Range("L34").Value = 272730120
$ws.Range("M34").Value = -1051.00002
$ws.Range("N34").Value = -272730288
$ws.Range("H69").Value = 962
$ws.Range("I69").Value = 831.3333
$ws.Range("J69").Value = 1550
$ws.Range("K69").Value = 2493.9999
$ws.Range("L69").Value = 4650
$ws.Range("M69").Value = -1682.9999
$ws.Range("N69").Value = -6272
$ws.Range("H72").Value = 962
$ws.Range("I72").Value = 831.3333
$ws.Range("J72").Value = 1550
$ws.Range("K72").Value = 7481.9997
$ws.Range("L72").Value = 13950
$ws.Range("M72").Value = -3425.9997
$ws.Range("N72").Value = -22062
$ws.Range("H135").Value = 932.7241
$ws.Range("I135").Value = 783.2222
$ws.Range("K135").Value = 7048.999800000001
$ws.Range("M135").Value = -4513.999800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 22763.334
$ws.Range("J32").Value = 22763.334
$ws.Range("L32").Value = 22763.334
$ws.Range("N32").Value = -23355.334
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").Value = $null
$ws.Range("H80").Value = 7165
$ws.Range("I80").Value = 2638.3333
$ws.Range("J80").Value = 11691.667
$ws.Range("K80").Value = 2638.3333
$ws.Range("L80").Value = 11691.667
$ws.Range("M80").Value = -1640.3333
$ws.Range("N80").Value = -13687.667
$ws.Range("H83").Value = 7165
$ws.Range("I83").Value = 2638.3333
$ws.Range("J83").Value = 11691.667
$ws.Range("K83").Value = 13191.6665
$ws.Range("L83").Value = 58458.335
$ws.Range("M83").Value = -8199.666499999999
$ws.Range("N83").Value = -68442.33499999999
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null
$ws.Range("H97").Value = 58824856
$ws.Range("I97").Value = 76924440
$ws.Range("J97").Value = 1190
$ws.Range("K97").Value = 76924440
$ws.Range("L97").Value = 1190
$ws.Range("M97").Value = -76923944
$ws.Range("N97").Value = -2182
$ws.Range("H102").Value = 1376.5
$ws.Range("I102").Value = 1235.4667
$ws.Range("J102").Value = 1799.6
$ws.Range("K102").Value = 1235.4667
$ws.Range("L102").Value = 1799.6
$ws.Range("M102").Value = 386.5333000000001
$ws.Range("N102").Value = -5043.6
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1575
$ws.Range("I46").Value = 1557.1428
$ws.Range("J46").Value = 1600
$ws.Range("K46").Value = 1557.1428
$ws.Range("L46").Value = 1600
$ws.Range("M46").Value = -1369.1428
$ws.Range("N46").Value = -1976
$ws.Range("H82").Value = 2333.25
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 2444.3333
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 2444.3333
$ws.Range("M82").Value = -1639
$ws.Range("N82").Value = -3166.3333
$ws.Range("H85").Value = 2333.25
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 2444.3333
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 2444.3333
$ws.Range("M85").Value = -752
$ws.Range("N85").Value = -4940.3333
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null
$ws.Range("H100").Value = 3700.0833
$ws.Range("I100").Value = 1980.6
$ws.Range("K100").Value = 1980.6
$ws.Range("M100").Value = -1439.6
$ws.Range("H136").Value = 2560.8696
$ws.Range("I136").Value = 1504.1666
$ws.Range("J136").Value = 3713.6365
$ws.Range("K136").Value = 4512.4998
$ws.Range("L136").Value = 11140.9095
$ws.Range("M136").Value = -1962.4998
$ws.Range("N136").Value = -16240.9095

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = $null
$ws.Range("H81").Value = 4222.773
$ws.Range("I81").Value = 1816.8334
$ws.Range("J81").Value = 5125
$ws.Range("K81").Value = 3633.6668
$ws.Range("L81").Value = 10250
$ws.Range("M81").Value = -2572.6668
$ws.Range("N81").Value = -12372
$ws.Range("H84").Value = 4222.773
$ws.Range("I84").Value = 1816.8334
$ws.Range("J84").Value = 5125
$ws.Range("K84").Value = 18168.334
$ws.Range("L84").Value = 51250
$ws.Range("M84").Value = -12864.334
$ws.Range("N84").Value = -61858
$ws.Range("H96").Value = 16687666
$ws.Range("I96").Value = 100000000
$ws.Range("J96").Value = 25199.4
$ws.Range("K96").Value = 100000000
$ws.Range("L96").Value = 25199.4
$ws.Range("M96").Value = -99998627
$ws.Range("N96").Value = -27945.4
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").Value = $null
$ws.Range("H127").Value = 44953.332
$ws.Range("J127").Value = 44953.332
$ws.Range("L127").Value = 44953.332
$ws.Range("N127").Value = -54873.332
